$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add a new client record in row 28, keeping the same layout as the
# existing rows (column D "Endereço" is left blank, as in all other rows).
# Numeric-looking values (password, CPF, CEP, phone) must stay as TEXT,
# exactly like every other row in this sheet, so force text format before
# assigning them - otherwise Excel would auto-coerce them to numbers.
$row = 28

$ws.Cells.Item($row, 1).Value = "francine benedetto"

$ws.Cells.Item($row, 2).NumberFormat = "@"
$ws.Cells.Item($row, 2).Value = "945833"

$ws.Cells.Item($row, 3).NumberFormat = "@"
$ws.Cells.Item($row, 3).Value = "85081450049"

$ws.Cells.Item($row, 5).NumberFormat = "@"
$ws.Cells.Item($row, 5).Value = "92410480"

$ws.Cells.Item($row, 6).Value = "francinebenedetto@gmail.com"

$ws.Cells.Item($row, 7).NumberFormat = "@"
$ws.Cells.Item($row, 7).Value = "51994723632"

$ws.Cells.Item($row, 8).Value = "rua tapajos 50"
